$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.461.01'
$ws.Range("E2").Value = '  -2.95%  '
$ws.Range("D3").Value = '2.463.99'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''311.30'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").Value = '''93.90'
$ws.Range("E6").Value = '  -6.77%  '
$ws.Range("D7").Value = '''0.551'
$ws.Range("E7").Value = '  -3.36%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '''0.507'
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("D10").Value = '''33.44'
$ws.Range("E10").Value = '  -7.19%  '
$ws.Range("D11").Value = '''0.0782'
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '''6.94'
$ws.Range("E13").Value = '  -5.62%  '
$ws.Range("D14").Value = '2.850.35'
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").Value = '2.484.71'
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("D16").Value = '''14.51'
$ws.Range("E16").Value = '  -8.74%  '
$ws.Range("D17").Value = '''0.789'
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").Value = '41.455.69'
$ws.Range("E18").Value = '  -2.93%  '
$ws.Range("D19").Value = '''6.35'
$ws.Range("E19").Value = '  -6.17%  '
$ws.Range("D20").Value = '0.0₃0916'
$ws.Range("E20").Value = '  -3.86%  '
$ws.Range("D21").Value = '''11.54'
$ws.Range("E21").Value = '  -5.68%  '
$ws.Range("D22").Value = '''68.64'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").Value = '''237.71'
$ws.Range("E23").Value = '  -2.50%  '
$ws.Range("D24").Value = '''2.77'
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("D25").Value = '''1.93'
$ws.Range("E25").Value = '  -5.48%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '''24.76'
$ws.Range("E27").Value = '  -5.01%  '
$ws.Range("D28").Value = '''2.21'
$ws.Range("E28").Value = '  -5.84%  '
$ws.Range("D29").Value = '''9.72'
$ws.Range("E29").Value = '  -4.30%  '
$ws.Range("D30").Value = '''36.33'
$ws.Range("E30").Value = '  -7.66%  '
$ws.Range("D31").Value = '''153.19'
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("D32").Value = '''5.64'
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = '''2.56'
$ws.Range("E34").Value = '  -7.31%  '
$ws.Range("D35").Value = '''0.0753'
$ws.Range("E35").Value = '  -5.07%  '
$ws.Range("D36").Value = '''3.02'
$ws.Range("E36").Value = '  -4.75%  '
$ws.Range("D37").Value = '''17.11'
$ws.Range("E37").Value = '  -6.77%  '
$ws.Range("D38").Value = '''1.88'
$ws.Range("E38").Value = '  -7.33%  '
$ws.Range("E39").Value = '  -6.16%  '
$ws.Range("E40").Value = '  -4.04%  '
$ws.Range("D41").Value = '''4.10'
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("D42").Value = '''21.33'
$ws.Range("E42").Value = '  -3.16%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '1.982.34'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '''0.0285'
$ws.Range("E45").Value = '  -4.63%  '
$ws.Range("D46").Value = '''3.05'
$ws.Range("E46").Value = '  -7.74%  '
$ws.Range("D47").Value = '''8.78'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("D48").Value = '''77.15'
$ws.Range("E48").Value = '  -5.05%  '
$ws.Range("D49").Value = '''97.50'
$ws.Range("E49").Value = '  -3.58%  '
$ws.Range("D50").Value = '''68.98'
$ws.Range("E50").Value = '  -4.99%  '
$ws.Range("D51").Value = '''0.180'
$ws.Range("E51").Value = '  -6.54%  '
